$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.517.70'
$ws.Range("E2").Value = '  +5.86%  '

# Row 3
$ws.Range("D3").Value = '2.740.59'
$ws.Range("E3").Value = '  +4.51%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.38'
$ws.Range("E5").Value = '  +6.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '332.34'
$ws.Range("E6").Value = '  +3.26%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.532'
$ws.Range("E7").Value = '  +2.33%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.572'
$ws.Range("E9").Value = '  +6.50%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.29'
$ws.Range("E10").Value = '  +5.47%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.21'
$ws.Range("E11").Value = '  +2.23%  '

# Row 12
$ws.Range("E12").Value = '  +2.94%  '

# Row 13
$ws.Range("E13").Value = '  +2.94%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.55'
$ws.Range("E14").Value = '  +5.09%  '

# Row 15
$ws.Range("D15").Value = '3.174.37'
$ws.Range("E15").Value = '  +4.55%  '

# Row 16
$ws.Range("D16").Value = '2.746.13'
$ws.Range("E16").Value = '  +3.66%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.881'
$ws.Range("E17").Value = '  +2.56%  '

# Row 18
$ws.Range("D18").Value = '51.436.47'
$ws.Range("E18").Value = '  +5.77%  '

# Row 19
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.40'
$ws.Range("E19").Value = '  +5.04%  '

# Row 20
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.07'
$ws.Range("E20").Value = '  +6.00%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.80'
$ws.Range("E21").Value = '  +2.13%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").Value = '  +2.17%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '277.23'
$ws.Range("E23").Value = '  +3.08%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.33'
$ws.Range("E24").Value = '  +1.03%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.63'
$ws.Range("E25").Value = '  +4.47%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.65'
$ws.Range("E26").Value = '  +2.64%  '

# Row 27
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  +1.41%  '

# Row 29
$ws.Range("E29").Value = '  -0.07%  '

# Row 30
$ws.Range("E30").Value = '  +1.67%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.92'
$ws.Range("E31").Value = '  +0.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.14'
$ws.Range("E32").Value = '  +1.73%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.54'
$ws.Range("E33").Value = '  +1.66%  '

# Row 34
$ws.Range("E34").Value = '  +3.30%  '

# Row 35
$ws.Range("E35").Value = '  +0.01%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.03'
$ws.Range("E36").Value = '  -0.98%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.07'
$ws.Range("E37").Value = '  +2.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.89'
$ws.Range("E38").Value = '  -0.64%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.14'
$ws.Range("E39").Value = '  +0.31%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '129.49'
$ws.Range("E40").Value = '  +3.32%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0347'
$ws.Range("E41").Value = '  +10.79%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.12'
$ws.Range("E42").Value = '  +3.20%  '

# Row 43
$ws.Range("E43").Value = '  +2.47%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.24'
$ws.Range("E44").Value = '  +4.48%  '

# Row 45
$ws.Range("E45").Value = '  +13.96%  '

# Row 46
$ws.Range("D46").Value = '2.101.76'
$ws.Range("E46").Value = '  +1.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("E47").Value = '  +3.20%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.22'
$ws.Range("E48").Value = '  +2.44%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.51'
$ws.Range("E49").Value = '  +7.21%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.91'
$ws.Range("E50").Value = '  -0.22%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.76'
$ws.Range("E51").Value = '  +2.27%  '
